# Replace the single "Slownie" amount-in-words run with the (accidentally
# captured) HTML error-page dump, preserving the run's original character
# formatting (rPr) and keeping it a separate run from its neighbours.
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute('dwadzieścia siedem tysięcy osiemset dziewięćdziesiąt sześć PLN 40/100', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "target run text not found"
}

$startPos = $rng.Start

# Nudge the formatting slightly before the text edit so this run is not
# silently coalesced into the identically-formatted ": " run right before
# it; we restore the exact original size immediately after.
$originalSize = $rng.Font.Size
$rng.Font.Size = $originalSize + 1

$lines = @(
    '<HTML>',
    '<HEAD>',
    '<META HTTP-EQUIV="Content-Type" CONTENT="text/html; charset=ISO-8859-2">',
    '<TITLE>',
    '(none)',
    '</TITLE>',
    '</HEAD>',
    '<BODY BGCOLOR="#FFFFFF">',
    '<H2 ALIGN="center">',
    '[Error : B³±d]',
    '</H2>',
    '<P>',
    '<H2 ALIGN="center">',
    '[403] Forbidden : Dostêp zabroniony',
    '</H2>',
    '<HR>',
    '<CENTER>',
    '<TABLE>',
    '<TR><TD VALIGN="top" ALIGN="left" WIDTH="50">',
    '<IMG SRC="/icons/gb.gif" WIDTH="26" HIGHT="19" ALT="English">',
    '</TD><TD>',
    'You don''t have permission to access',
    '<STRONG>',
    '/slownie.php?format=3&amp;kwota=27896.40',
    '</STRONG>',
    'on this server.',
    '</TD></TR><TR><TD></TD><TD><HR></TD></TR>',
    '<TR><TD VALIGN="top" ALIGN="left" WIDTH="50">',
    '<IMG SRC="/icons/poland.gif" WIDTH="26" HIGHT="19" ALT="Polish">',
    '</TD><TD>',
    'Zapytanie odrzucone przez serwer. Nie masz dostêpu do',
    '<STRONG>',
    '/slownie.php?format=3&amp;kwota=27896.40',
    '</STRONG>',
    'na tym serwerze.',
    '</TD></TR>',
    '</TABLE>',
    '</CENTER>',
    '<P>',
    '<HR>',
    '<P>',
    '<CENTER>',
    '<IMG SRC="/icons/email.gif" WIDTH="50" HIGHT="66" ALT="Mail to:">',
    '<BR>',
    '<A HREF="/cdn-cgi/l/email-protection#9fe8fafdf2feecebfaeddffeeff6b1f0f9f9f6fcfafdf3f0f8b1eff3a0eceafdf5fafceba2ddf3fefba5bfc4abafacc2bf959696d9f0edfdf6fbfbfaf1b3bfcacdd3a2b0ecf3f0e8f1f6fab1eff7efa0f9f0edf2feeba2acb9fef2efa4f4e8f0ebfea2ada8a7a6a9b1abaf">',
    '<span class="__cf_email__" data-cfemail="2a5d4f48474b595e4f586a4b5a4304454c4c43494f4846454d045a46">[email&#160;protected]</span>',
    '</A>',
    '<P>',
    '<TABLE WIDTH="50%">',
    '<TR>',
    '<TD ALIGN="left">',
    '<IMG SRC="/icons/apache_pb.gif" WIDTH="259" HIGHT="32" ALT="Powered by Apache">',
    '</TD>',
    '<TD ALIGN="right">',
    '<IMG SRC="/icons/linux_pwd.gif" WIDTH="196" HIGHT="49" ALT="Powered by Linux">',
    '</TD>',
    '</TR>',
    '</TABLE>',
    '<BR>',
    '<FONT SIZE="-3">',
    'data:Tuesday, 16-Apr-2024 20:37:28 CEST,',
    'ostatnia modyfikacja: Tuesday, 02-Nov-2021 12:49:14 CET',
    '</FONT>',
    '</CENTER>',
    '<script data-cfasync="false" src="/cdn-cgi/scripts/5c5dd728/cloudflare-static/email-decode.min.js"></script></BODY>',
    '</HTML>'
)

$lineBreak = [char]11
$newText = ($lines -join $lineBreak) + $lineBreak

$rng.Text = $newText

$endPos = $startPos + $newText.Length
$rng2 = $d.Range($startPos, $endPos)
$rng2.Font.Size = $originalSize

Write-Output "replaced run with"
Write-Output $lines.Count
